$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $CellRef, $Text) {
    $Sheet.Range($CellRef).Value = "'" + $Text
    $Sheet.Range($CellRef).Style = "Normal"
}

Set-TextValue $ws "D2" '258.28'
Set-TextValue $ws "E2" '0.18%'
Set-TextValue $ws "D3" '27.07'
Set-TextValue $ws "E3" '-0.51%'
Set-TextValue $ws "D4" '4.608'
Set-TextValue $ws "E4" '-5.51%'
Set-TextValue $ws "D5" '0.05908'
Set-TextValue $ws "E5" '-0.64%'
Set-TextValue $ws "D6" '6.644'
Set-TextValue $ws "E6" '-0.81%'
Set-TextValue $ws "D7" '0.8541'
Set-TextValue $ws "E7" '-1.66%'
Set-TextValue $ws "D8" '0.9454'
Set-TextValue $ws "E8" '-5.30%'
Set-TextValue $ws "D9" '0.1404'
Set-TextValue $ws "E9" '-0.56%'
Set-TextValue $ws "D10" '0.05330'
Set-TextValue $ws "E10" '50.01%'
Set-TextValue $ws "D11" '0.07098'
Set-TextValue $ws "E11" '-1.18%'
Set-TextValue $ws "D12" '0.03114'
Set-TextValue $ws "E12" '-1.11%'
Set-TextValue $ws "D13" '0.09138'
Set-TextValue $ws "E13" '-1.23%'
Set-TextValue $ws "D14" '0.001536'
Set-TextValue $ws "E14" '-0.20%'
Set-TextValue $ws "B15" 'TigerCash'
Set-TextValue $ws "C15" 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws "D15" '0.006009'
Set-TextValue $ws "E15" '-0.55%'
Set-TextValue $ws "B16" 'LEO'
Set-TextValue $ws "C16" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws "D16" '3.506'
Set-TextValue $ws "E16" '0.44%'
Set-TextValue $ws "B17" 'GateToken'
Set-TextValue $ws "C17" 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws "D17" '3.184'
Set-TextValue $ws "E17" '-2.39%'
Set-TextValue $ws "B18" 'BTSEToken'
Set-TextValue $ws "C18" 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws "D18" '2.204'
Set-TextValue $ws "E18" '0.00%'
Set-TextValue $ws "B19" 'One'
Set-TextValue $ws "C19" 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws "D19" '0.0006060'
Set-TextValue $ws "E19" '-0.33%'
Set-TextValue $ws "D20" '0.3056'
Set-TextValue $ws "E20" '-2.91%'
Set-TextValue $ws "E21" '-2.18%'
Set-TextValue $ws "D22" '3.818'
Set-TextValue $ws "E22" '8.22%'
Set-TextValue $ws "E23" '0.12%'
Set-TextValue $ws "D24" '0.001221'
Set-TextValue $ws "E24" '0.00%'
Set-TextValue $ws "D25" '0.004289'
Set-TextValue $ws "E25" '-5.06%'
Set-TextValue $ws "E26" '-0.01%'
Set-TextValue $ws "D27" '0.0001937'
Set-TextValue $ws "E27" '29.91%'
Set-TextValue $ws "D40" '0.03827'
Set-TextValue $ws "E40" '-0.11%'
Set-TextValue $ws "D41" '0.006250'
Set-TextValue $ws "E41" '57.46%'
Set-TextValue $ws "D42" '0.1101'
Set-TextValue $ws "E42" '-0.24%'
Set-TextValue $ws "D43" '0.002200'
Set-TextValue $ws "E43" '-4.77%'
Set-TextValue $ws "D44" '0.01238'
Set-TextValue $ws "E44" '17.98%'
Set-TextValue $ws "D45" '0.00005469'
Set-TextValue $ws "E45" '-0.46%'
Set-TextValue $ws "E46" '0.01%'
Set-TextValue $ws "D47" '0.05100'
Set-TextValue $ws "D48" '0.2502'
Set-TextValue $ws "E48" '11,519.85%'
Set-TextValue $ws "E49" '0.01%'
Set-TextValue $ws "E50" '0.01%'
